$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Update the date line at the top of the document
$d.Content.Find.Execute("2025-10-20 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-21 Tuesday", 2)

# Row 1
$t.Cell(1, 1).Range.Text = "16+29=45"
$t.Cell(1, 2).Range.Text = "76-28=48"
$t.Cell(1, 3).Range.Text = "36+26=62"
$t.Cell(1, 4).Range.Text = "26+25=51"
$t.Cell(1, 5).Range.Text = "24+19=43"

# Row 2
$t.Cell(2, 1).Range.Text = "36+58=94"
$t.Cell(2, 2).Range.Text = "71-56=15"
$t.Cell(2, 3).Range.Text = "48+33=81"
$t.Cell(2, 4).Range.Text = "19+6=25"
$t.Cell(2, 5).Range.Text = "9+34=43"

# Row 3
$t.Cell(3, 1).Range.Text = "27+19=46"
$t.Cell(3, 2).Range.Text = "39+27=66"
$t.Cell(3, 3).Range.Text = "6+67=73"
$t.Cell(3, 4).Range.Text = "46+25=71"
$t.Cell(3, 5).Range.Text = "38+45=83"

# Row 4
$t.Cell(4, 1).Range.Text = "5+39=44"
$t.Cell(4, 2).Range.Text = "20-12=8"
$t.Cell(4, 3).Range.Text = "57-18=39"
$t.Cell(4, 4).Range.Text = "40-13=27"
$t.Cell(4, 5).Range.Text = "27+54=81"

# Row 5
$t.Cell(5, 1).Range.Text = "61-42=19"
$t.Cell(5, 2).Range.Text = "29+52=81"
$t.Cell(5, 3).Range.Text = "58+4=62"
$t.Cell(5, 4).Range.Text = "84-39=45"
$t.Cell(5, 5).Range.Text = "17+69=86"

# Row 6
$t.Cell(6, 1).Range.Text = "63-6=57"
$t.Cell(6, 2).Range.Text = "39+24=63"
$t.Cell(6, 3).Range.Text = "15+36=51"
$t.Cell(6, 4).Range.Text = "94-35=59"
$t.Cell(6, 5).Range.Text = "91-83=8"

# Row 7
$t.Cell(7, 1).Range.Text = "36+39=75"
$t.Cell(7, 2).Range.Text = "17+28=45"
$t.Cell(7, 3).Range.Text = "82-68=14"
$t.Cell(7, 4).Range.Text = "94-36=58"
$t.Cell(7, 5).Range.Text = "83-35=48"

# Row 8
$t.Cell(8, 1).Range.Text = "69+7=76"
$t.Cell(8, 2).Range.Text = "22+49=71"
$t.Cell(8, 3).Range.Text = "27+7=34"
$t.Cell(8, 4).Range.Text = "91-53=38"
$t.Cell(8, 5).Range.Text = "72-53=19"

# Row 9
$t.Cell(9, 1).Range.Text = "7+16=23"
$t.Cell(9, 2).Range.Text = "56-27=29"
$t.Cell(9, 3).Range.Text = "16+48=64"
$t.Cell(9, 4).Range.Text = "35+9=44"
$t.Cell(9, 5).Range.Text = "60-38=22"

# Row 10
$t.Cell(10, 1).Range.Text = "84+9=93"
$t.Cell(10, 2).Range.Text = "37+18=55"
$t.Cell(10, 3).Range.Text = "26+29=55"
$t.Cell(10, 4).Range.Text = "65+9=74"
$t.Cell(10, 5).Range.Text = "42-15=27"

# Row 11
$t.Cell(11, 1).Range.Text = "42-39=3"
$t.Cell(11, 2).Range.Text = "17+6=23"
$t.Cell(11, 3).Range.Text = "26+66=92"
$t.Cell(11, 4).Range.Text = "44-5=39"
$t.Cell(11, 5).Range.Text = "75-38=37"

# Row 12
$t.Cell(12, 1).Range.Text = "18+54=72"
$t.Cell(12, 2).Range.Text = "13-9=4"
$t.Cell(12, 3).Range.Text = "33+8=41"
$t.Cell(12, 4).Range.Text = "41-12=29"
$t.Cell(12, 5).Range.Text = "35+7=42"

# Row 13
$t.Cell(13, 1).Range.Text = "60-37=23"
$t.Cell(13, 2).Range.Text = "97-38=59"
$t.Cell(13, 3).Range.Text = "33+48=81"
$t.Cell(13, 4).Range.Text = "36+17=53"
$t.Cell(13, 5).Range.Text = "46-7=39"

# Row 14
$t.Cell(14, 1).Range.Text = "47+37=84"
$t.Cell(14, 2).Range.Text = "76+18=94"
$t.Cell(14, 3).Range.Text = "28+47=75"
$t.Cell(14, 4).Range.Text = "37+34=71"
$t.Cell(14, 5).Range.Text = "67+15=82"

# Row 15
$t.Cell(15, 1).Range.Text = "18+74=92"
$t.Cell(15, 2).Range.Text = "4+59=63"
$t.Cell(15, 3).Range.Text = "30-19=11"
$t.Cell(15, 4).Range.Text = "78-9=69"
$t.Cell(15, 5).Range.Text = "95-48=47"

# Row 16
$t.Cell(16, 1).Range.Text = "52-29=23"
$t.Cell(16, 2).Range.Text = "18+66=84"
$t.Cell(16, 3).Range.Text = "8+15=23"
$t.Cell(16, 4).Range.Text = "67-29=38"
$t.Cell(16, 5).Range.Text = "50-44=6"

# Row 17
$t.Cell(17, 1).Range.Text = "91-59=32"
$t.Cell(17, 2).Range.Text = "44+29=73"
$t.Cell(17, 3).Range.Text = "64-39=25"
$t.Cell(17, 4).Range.Text = "58-9=49"
$t.Cell(17, 5).Range.Text = "43-19=24"

# Row 18
$t.Cell(18, 1).Range.Text = "91-74=17"
$t.Cell(18, 2).Range.Text = "45+48=93"
$t.Cell(18, 3).Range.Text = "54-46=8"
$t.Cell(18, 4).Range.Text = "73-9=64"
$t.Cell(18, 5).Range.Text = "96-88=8"

# Row 19
$t.Cell(19, 1).Range.Text = "51-49=2"
$t.Cell(19, 2).Range.Text = "91-23=68"
$t.Cell(19, 3).Range.Text = "75-37=38"
$t.Cell(19, 4).Range.Text = "60-29=31"
$t.Cell(19, 5).Range.Text = "5+86=91"

# Row 20
$t.Cell(20, 1).Range.Text = "86-37=49"
$t.Cell(20, 2).Range.Text = "94-25=69"
$t.Cell(20, 3).Range.Text = "80-73=7"
$t.Cell(20, 4).Range.Text = "68+7=75"
$t.Cell(20, 5).Range.Text = "58+35=93"
